# Rename Sheet3 -> addNewCustomer
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "addNewCustomer"

# Extend the bordered/formatted block from A1:A5 to A1:I5, reusing the
# existing cell style (border) rather than creating new style entries.
$ws.Range("A1:A5").Copy()
$ws.Range("A1:I5").PasteSpecial(-4122)   # xlPasteFormats

# Cell values, entered in the same order as the original authoring session
# (this keeps the sharedStrings table ordering identical to the recording).
$ws.Range("A1").Value = "companyName"
$ws.Range("A3").Value = "Viettel"
$ws.Range("B1").Value = "vatNumber"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "website"
$ws.Range("E1").Value = "groups"
$ws.Range("F1").Value = "address"
$ws.Range("G1").Value = "city"
$ws.Range("H1").Value = "state"
$ws.Range("I1").Value = "zipCode"

$ws.Range("A4").Value = "HBT"

# B4 / C4 are numeric-looking text, entered with a leading apostrophe so
# they are stored as text (quote-prefixed) instead of numbers - this keeps
# the leading zero in the phone number, matching the recorded workbook.
$ws.Range("B4").Value = "'667735"
$ws.Range("C4").Value = "'0852741963"

$ws.Range("D4").Value = "hbt.hn.com"
$ws.Range("E4").Value = "khoinghia"
$ws.Range("F4").Value = "Me Linh"
$ws.Range("G4").Value = "Ha Noi"
$ws.Range("H4").Value = "Pass"
$ws.Range("I4").Value = 28386

# Column widths (best effort - COM rounds to the nearest pixel column width)
$ws.Columns.Item(1).ColumnWidth = 20.6
$ws.Columns.Item(2).ColumnWidth = 12.6
$ws.Columns.Item(3).ColumnWidth = 13.1
$ws.Columns.Item(4).ColumnWidth = 11.6
$ws.Range("E1:F1").ColumnWidth = 10.0

# Selection moves to A5
$ws.Range("A5").Select() | Out-Null
